$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 678
$ws.Range("F6").Value = 261
$ws.Range("F9").Value = 7350
$ws.Range("F17").Value = 1790
$ws.Range("F18").Value = 1084
$ws.Range("F20").Value = 72
$ws.Range("F21").Value = 1833
$ws.Range("F22").Value = 1375
$ws.Range("F25").Value = 54
$ws.Range("F26").Value = 1131
$ws.Range("F27").Value = 121
$ws.Range("F28").Value = 537
$ws.Range("F30").Value = 72
$ws.Range("F31").Value = 3916
$ws.Range("F33").Value = 3894
$ws.Range("F36").Value = 224
$ws.Range("F39").Value = 47
$ws.Range("F41").Value = 379
$ws.Range("F44").Value = 254
$ws.Range("F46").Value = 772
$ws.Range("F48").Value = 3
$ws.Range("F49").Value = 140

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 3
$ws.Range("F17").Value = 548
$ws.Range("F24").Value = 89

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1849
$ws.Range("F8").Value = 2896
$ws.Range("F9").Value = 1135
$ws.Range("F10").Value = 1124
$ws.Range("F12").Value = 429
$ws.Range("F13").Value = 1838
$ws.Range("F14").Value = 8186

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 678
$ws.Range("F5").Value = 261
$ws.Range("F6").Value = 1849
$ws.Range("F8").Value = 2896
$ws.Range("F9").Value = 7350
$ws.Range("F10").Value = 1135
$ws.Range("F11").Value = 1124
$ws.Range("F13").Value = 429
$ws.Range("F18").Value = 1084
$ws.Range("F20").Value = 72
$ws.Range("F21").Value = 1833
$ws.Range("F22").Value = 1375
$ws.Range("F25").Value = 54
$ws.Range("F26").Value = 1131
$ws.Range("F28").Value = 121
$ws.Range("F30").Value = 548
$ws.Range("F31").Value = 537
$ws.Range("F34").Value = 72
$ws.Range("F35").Value = 3917
$ws.Range("F37").Value = 3894
$ws.Range("F40").Value = 224
$ws.Range("F44").Value = 379
$ws.Range("F45").Value = 89
$ws.Range("F47").Value = 254
